# Update the date header
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-12-25 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-26 Thursday", 2)

# Update the math problems table.
# The table has 20 rows; data rows are 1, 5, 9, 13, 17 (1-based), each with 5 cells.
# Replacing cell text directly (in row/column order) avoids collisions between
# old and new values that would otherwise occur with a global find & replace.
$t = $d.Tables.Item(1)

$rowsData = @{
    1  = @("82÷6=", "29÷3=", "72÷6=", "22÷3=", "55÷8=")
    5  = @("76÷4=", "62÷5=", "93÷5=", "95÷3=", "36÷8=")
    9  = @("42÷6=", "38÷5=", "22÷2=", "95÷8=", "18÷2=")
    13 = @("27÷5=", "46÷8=", "78÷7=", "80÷3=", "19÷7=")
    17 = @("82÷3=", "12÷6=", "40÷4=", "57÷9=", "29÷2=")
}

foreach ($rowIndex in $rowsData.Keys) {
    $values = $rowsData[$rowIndex]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
